$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 16666.666
$ws.Range("I69").Value = 15000
$ws.Range("J69").Value = 20000
$ws.Range("K69").Value = 45000
$ws.Range("L69").Value = 60000
$ws.Range("M69").Value = -44126
$ws.Range("N69").Value = -61748

$ws.Range("H72").Value = 16666.666
$ws.Range("I72").Value = 15000
$ws.Range("J72").Value = 20000
$ws.Range("K72").Value = 135000
$ws.Range("L72").Value = 180000
$ws.Range("M72").Value = -130632
$ws.Range("N72").Value = -188736

$ws.Range("H137").Value = 2640603.2
$ws.Range("I137").Value = 3847766
$ws.Range("K137").Value = 11543298
$ws.Range("M137").Value = -11540748

$ws.Range("H138").Value = 6637.478
$ws.Range("J138").Value = 7340.4736
$ws.Range("L138").Value = 22021.4208
$ws.Range("N138").Value = -32301.4208

$ws.Range("H141").Value = 2741.5
$ws.Range("I141").Value = 2741.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8224.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -3044.5
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1265126.6
$ws.Range("I32").Value = 618940.9
$ws.Range("K32").Value = 618940.9
$ws.Range("M32").Value = -618653.9

$ws.Range("H61").Value = 2563.3845
$ws.Range("I61").Value = 2243.8948
$ws.Range("K61").Value = 2243.8948
$ws.Range("M61").Value = -2031.8948

$ws.Range("H74").Value = 149036.89
$ws.Range("I74").Value = 224443.28
$ws.Range("K74").Value = 224443.28
$ws.Range("M74").Value = -223569.28

$ws.Range("H77").Value = 149036.89
$ws.Range("I77").Value = 224443.28
$ws.Range("K77").Value = 1122216.4
$ws.Range("M77").Value = -1117848.4

$ws.Range("H122").Value = 2616.4102
$ws.Range("I122").Value = 2487.5652
$ws.Range("J122").Value = 2801.625
$ws.Range("K122").Value = 7462.6956
$ws.Range("L122").Value = 8404.875
$ws.Range("M122").Value = -5012.6956
$ws.Range("N122").Value = -13304.875

$ws.Range("H132").Value = 3665530
$ws.Range("I132").Value = 5496652
$ws.Range("K132").Value = 16489956
$ws.Range("M132").Value = -16487426

$ws.Range("H136").Value = 2563.3845
$ws.Range("I136").Value = 2243.8948
$ws.Range("K136").Value = 6731.6844
$ws.Range("M136").Value = -4181.6844

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4974.1665
$ws.Range("I99").Value = 4469.1
$ws.Range("K99").Value = 4469.1
$ws.Range("M99").Value = -2971.1

$ws.Range("H105").Value = 8966646
$ws.Range("I105").Value = 417544.22
$ws.Range("K105").Value = 417544.22
$ws.Range("M105").Value = -415797.22

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 620.619
$ws.Range("I7").Value = 767
$ws.Range("J7").Value = 152.2
$ws.Range("K7").Value = 767
$ws.Range("L7").Value = 152.2
$ws.Range("M7").Value = -654
$ws.Range("N7").Value = -378.2

$ws.Range("H31").Value = 3292666.2
$ws.Range("I31").Value = 1668.8572
$ws.Range("J31").Value = 5212414.5
$ws.Range("K31").Value = 1668.8572
$ws.Range("L31").Value = 5212414.5
$ws.Range("M31").Value = -1373.8572
$ws.Range("N31").Value = -5213004.5

$ws.Range("H34").Value = 3292666.2
$ws.Range("I34").Value = 1668.8572
$ws.Range("J34").Value = 5212414.5
$ws.Range("K34").Value = 1668.8572
$ws.Range("L34").Value = 5212414.5
$ws.Range("M34").Value = -1466.8572
$ws.Range("N34").Value = -5212818.5

$ws.Range("H58").Value = 4116.6665
$ws.Range("I58").Value = 350
$ws.Range("K58").Value = 350
$ws.Range("M58").Value = -147

$ws.Range("H62").Value = 3791.0833
$ws.Range("I62").Value = 3713.2856
$ws.Range("J62").Value = 3900
$ws.Range("K62").Value = 3713.2856
$ws.Range("L62").Value = 3900
$ws.Range("M62").Value = -3089.2856
$ws.Range("N62").Value = -5148

$ws.Range("H65").Value = 3791.0833
$ws.Range("I65").Value = 3713.2856
$ws.Range("J65").Value = 3900
$ws.Range("K65").Value = 18566.428
$ws.Range("L65").Value = 19500
$ws.Range("M65").Value = -15446.428
$ws.Range("N65").Value = -25740

$ws.Range("H86").Value = 7362.25
$ws.Range("I86").Value = 1997.6666
$ws.Range("J86").Value = 23456
$ws.Range("K86").Value = 1997.6666
$ws.Range("L86").Value = 23456
$ws.Range("M86").Value = -874.6666
$ws.Range("N86").Value = -25702

$ws.Range("H89").Value = 7362.25
$ws.Range("I89").Value = 1997.6666
$ws.Range("J89").Value = 23456
$ws.Range("K89").Value = 9988.333000000001
$ws.Range("L89").Value = 117280
$ws.Range("M89").Value = -4372.333000000001
$ws.Range("N89").Value = -128512

$ws.Range("H136").Value = 4116.6665
$ws.Range("I136").Value = 350
$ws.Range("K136").Value = 1050
$ws.Range("M136").Value = 1500

$ws.Range("H141").Value = 133442.89
$ws.Range("J141").Value = 141569.42
$ws.Range("L141").Value = 141569.42
$ws.Range("N141").Value = -151929.42

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 311
$ws.Range("I7").Value = 311
$ws.Range("K7").Value = 933
$ws.Range("M7").Value = -821

$ws.Range("H36").Value = 958.25
$ws.Range("I36").Value = 777.6667
$ws.Range("K36").Value = 2333.0001
$ws.Range("M36").Value = -2164.0001

$ws.Range("H56").Value = 7821.8887
$ws.Range("I56").Value = 7821.8887
$ws.Range("K56").Value = 7821.8887
$ws.Range("M56").Value = -7291.8887

$ws.Range("H137").Value = 7704.778
$ws.Range("I137").Value = 9306.143
$ws.Range("J137").Value = 2100
$ws.Range("K137").Value = 27918.429
$ws.Range("L137").Value = 6300
$ws.Range("M137").Value = -22818.429
$ws.Range("N137").Value = -16500

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 5605.8
$ws.Range("I2").Value = 61.875
$ws.Range("K2").Value = 61.875
$ws.Range("M2").Value = 51.125

$ws.Range("H132").Value = 2993.2222
$ws.Range("I132").Value = 3125.2666
$ws.Range("K132").Value = 9375.799800000001
$ws.Range("M132").Value = -6845.799800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 63391.8
$ws.Range("I40").Value = 73431.586
$ws.Range("K40").Value = 73431.586
$ws.Range("M40").Value = -73295.586

$ws.Range("H122").Value = 6608
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 6608
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 19824
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -24724

$ws.Range("H132").Value = 7846.353
$ws.Range("I132").Value = 9423.75
$ws.Range("K132").Value = 28271.25
$ws.Range("M132").Value = -25741.25

$ws.Range("H136").Value = 6580.4116
$ws.Range("I136").Value = 4127.143
$ws.Range("J136").Value = 8297.700000000001
$ws.Range("K136").Value = 12381.429
$ws.Range("L136").Value = 24893.1
$ws.Range("M136").Value = -9831.429
$ws.Range("N136").Value = -29993.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5268.4814
$ws.Range("I81").Value = 4589.357
$ws.Range("J81").Value = 5999.846
$ws.Range("K81").Value = 9178.714
$ws.Range("L81").Value = 11999.692
$ws.Range("M81").Value = -8117.714
$ws.Range("N81").Value = -14121.692

$ws.Range("H84").Value = 5268.4814
$ws.Range("I84").Value = 4589.357
$ws.Range("J84").Value = 5999.846
$ws.Range("K84").Value = 45893.57
$ws.Range("L84").Value = 59998.45999999999
$ws.Range("M84").Value = -40589.57
$ws.Range("N84").Value = -70606.45999999999

$ws.Range("H100").Value = 35715436
$ws.Range("I100").Value = 910.55
$ws.Range("J100").Value = 125001750
$ws.Range("K100").Value = 1821.1
$ws.Range("L100").Value = 250003500
$ws.Range("M100").Value = -1280.1
$ws.Range("N100").Value = -250004582

$ws.Range("H132").Value = 2792.0232
$ws.Range("I132").Value = 2647.675
$ws.Range("J132").Value = 4716.6665
$ws.Range("K132").Value = 7943.025000000001
$ws.Range("L132").Value = 14149.9995
$ws.Range("M132").Value = -5413.025000000001
$ws.Range("N132").Value = -19209.9995

$ws.Range("H136").Value = 10915.41
$ws.Range("I136").Value = 11077.6455
$ws.Range("J136").Value = 10286.75
$ws.Range("K136").Value = 33232.9365
$ws.Range("L136").Value = 30860.25
$ws.Range("M136").Value = -30682.9365
$ws.Range("N136").Value = -35960.25
